# Applies the commit's edits to the document:
#  - Collapses runs that were only split apart by (now-removed) proofErr
#    spell/grammar markers back into single runs, for several paragraphs.
#  - Inserts "Augusto" into the first "... persona ():" line and relocates
#    the "_GoBack" bookmark there (it previously sat in the "5ª persona(...)"
#    line).

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Paragraph 1 : "Tareas a realizar" / ":" -> single run "Tareas a
#    realizar:". The leading <w:proofErr w:type="gramStart"/> sits right
#    at the start of the paragraph (before any run), so a plain text
#    replace never reaches/clears it. Rebuild the paragraph instead: add
#    a fresh paragraph (inherits the same paragraph mark formatting),
#    fill it in, then delete the old paragraph (which carries the stray
#    proofErr) away.
# ---------------------------------------------------------------------
$p1 = $d.Paragraphs.Item(1)
$p1.Range.InsertParagraphAfter()
$newP1 = $d.Paragraphs.Item(2)
$newP1.Range.Text = "Tareas a realizar:"
$d.Paragraphs.Item(1).Range.Delete()

# ---------------------------------------------------------------------
# 2) Plain run-merges: replacing the full (multi-run) text with itself
#    makes the engine coalesce the runs into one and drop any proofErr
#    markers that were fully inside the replaced span.
# ---------------------------------------------------------------------
$merges = @(
    "1.Equipo de trabajo(30%):",
    "1.1 Configuración del equipo de trabajo(24%): ",
    "1.2. Nº de técnicos del equipo de trabajo(38%): ",
    "2.Proyecto técnico(30%):",
    "2.1.Solución propuesta en su conjunto(40%):",
    "Puntos: 1.2 , 1.3.",
    "2ª persona():",
    "3ª persona( ):",
    "4ª persona():"
)

foreach ($needle in $merges) {
    $find = $d.Content.Find
    $find.Execute($needle, $false, $false, $false, $false, $false, $true, 1, $false, $needle, 2) | Out-Null
}

# ---------------------------------------------------------------------
# 3) "1ª persona ():" -> "1ª persona (Augusto):" and move the "_GoBack"
#    bookmark from the "5ª persona(José Antonio):" paragraph into this
#    one, sitting right after "Augusto".
# ---------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text.StartsWith("1ª persona (")) {
        $target = $cand
        break
    }
}

$pStart = $target.Range.Start
$text = $target.Range.Text
$idx = $text.IndexOf("(")
$leftPos = $pStart + $idx + 1

$rIns = $d.Range($leftPos, $leftPos)
$rIns.InsertAfter("Augusto")
$rightPos = $leftPos + 7

# Bookmark at the left edge of "Augusto" forces the run split between
# "1ª persona (" and "Augusto"; it is removed again once that split has
# happened, leaving only the real "_GoBack" bookmark on the right edge.
$rLeft = $d.Range($leftPos, $leftPos)
$d.Bookmarks.Add("TempSplitMark", $rLeft)

$rRight = $d.Range($rightPos, $rightPos)
$d.Bookmarks.Add("_GoBack", $rRight)

$d.Bookmarks.Item("TempSplitMark").Delete()

Write-Output "done"
